$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump volume/issue number and report week date range ---
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Crime Complaints table (rows 15-28): refreshed weekly counts ---

# Row 15
$ws.Range("N15").Value = -40

# Row 16
$ws.Range("F16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("F16").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("H16").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("H16").Value = -50
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = -9.523809523809
$ws.Range("L16").Value = -38.709677419354
$ws.Range("M16").Value = -26.923076923076
$ws.Range("N16").Value = -84.166666666666

# Row 17
$ws.Range("C22").Copy($ws.Range("C17"))
$ws.Range("D22").Copy($ws.Range("D17"))
$ws.Range("E22").Copy($ws.Range("E17"))
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = -14.285714285714
$ws.Range("L17").Value = 4.166666666666
$ws.Range("M17").Value = 31.578947368421
$ws.Range("N17").Value = -37.5

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 36.842105263157
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -23.423423423423
$ws.Range("L18").Value = -15
$ws.Range("N18").Value = -73.4375

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -23.529411764705
$ws.Range("I19").Value = 151
$ws.Range("J19").Value = 236
$ws.Range("K19").Value = -36.016949152542
$ws.Range("L19").Value = -12.209302325581
$ws.Range("M19").Value = 17.96875
$ws.Range("N19").Value = -3.821656050955

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -25
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 120
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = 73.469387755102
$ws.Range("L20").Value = 157.575757575758
$ws.Range("M20").Value = 88.888888888888
$ws.Range("N20").Value = -91.842610364683

# Row 21
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -20.833333333333
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = 3.296703296703
$ws.Range("I21").Value = 368
$ws.Range("J21").Value = 449
$ws.Range("K21").Value = -18.040089086859
$ws.Range("L21").Value = 2.222222222222
$ws.Range("M21").Value = 20.655737704918
$ws.Range("N21").Value = -78.160237388724

# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 57.142857142857
$ws.Range("F24").Value = 33
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = -17.5
$ws.Range("I24").Value = 178
$ws.Range("J24").Value = 173
$ws.Range("K24").Value = 2.890173410404
$ws.Range("L24").Value = -31.007751937984
$ws.Range("M24").Value = 9.876543209876

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = 8.333333333333
$ws.Range("L25").Value = 14.705882352941

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 11
$ws.Range("H26").Value = 36.363636363636
$ws.Range("I26").Value = 59
$ws.Range("J26").Value = 65
$ws.Range("K26").Value = -9.230769230769
$ws.Range("L26").Value = -10.60606060606
$ws.Range("M26").Value = 28.260869565217

# Row 27
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50

# Row 28
$ws.Range("D22").Copy($ws.Range("D28"))
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
